$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 531.9
$ws.Range("C2").Value = 526
$ws.Range("D2").Value = 530.05
$ws.Range("E2").Value = 529.2
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 528.35

$ws.Range("B3").Value = 3190
$ws.Range("C3").Value = 3116
$ws.Range("D3").Value = 3129.95
$ws.Range("E3").Value = 3137.2
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 3186.8

$ws.Range("B4").Value = 510.95
$ws.Range("C4").Value = 501.25
$ws.Range("D4").Value = 504.7
$ws.Range("E4").Value = 505.05
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 505.6

$ws.Range("B5").Value = 1879.9
$ws.Range("C5").Value = 1855.1
$ws.Range("D5").Value = 1877.95
$ws.Range("E5").Value = 1876.4
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1864.3

$ws.Range("B6").Value = 7317.95
$ws.Range("C6").Value = 7226.05
$ws.Range("D6").Value = 7300.1
$ws.Range("E6").Value = 7302
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 7262.7

$ws.Range("B7").Value = 211.8
$ws.Range("C7").Value = 200.25
$ws.Range("D7").Value = 209.44
$ws.Range("E7").Value = 210.26
$ws.Range("F7").Value = 1298
$ws.Range("G7").Value = 201.35

$ws.Range("B8").Value = 248.2
$ws.Range("C8").Value = 242.13
$ws.Range("D8").Value = 242.48
$ws.Range("E8").Value = 242.42
$ws.Range("F8").Value = 94
$ws.Range("G8").Value = 246.75

$ws.Range("B9").Value = 493.9
$ws.Range("C9").Value = 488
$ws.Range("D9").Value = 493.1
$ws.Range("E9").Value = 492.95
$ws.Range("F9").Value = 43
$ws.Range("G9").Value = 488.6

$ws.Range("B10").Value = 859.75
$ws.Range("C10").Value = 841.75
$ws.Range("D10").Value = 845.8
$ws.Range("E10").Value = 846.6
$ws.Range("F10").Value = 21
$ws.Range("G10").Value = 854.35

$ws.Range("B11").Value = 4808.15
$ws.Range("C11").Value = 4664.15
$ws.Range("D11").Value = 4735.9
$ws.Range("E11").Value = 4728.05
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 4670.05

$ws.Range("B12").Value = 189.3
$ws.Range("C12").Value = 186.63
$ws.Range("D12").Value = 187.58
$ws.Range("E12").Value = 187.55
$ws.Range("F12").Value = 52
$ws.Range("G12").Value = 187.15

$ws.Range("B13").Value = 1852.8
$ws.Range("C13").Value = 1821
$ws.Range("D13").Value = 1839.95
$ws.Range("E13").Value = 1839.65
$ws.Range("F13").Value = 29
$ws.Range("G13").Value = 1845

$ws.Range("B14").Value = 1659.05
$ws.Range("C14").Value = 1643.25
$ws.Range("D14").Value = 1651.5
$ws.Range("E14").Value = 1651
$ws.Range("F14").Value = 118
$ws.Range("G14").Value = 1652.15

$ws.Range("B15").Value = 749.95
$ws.Range("C15").Value = 736.05
$ws.Range("D15").Value = 746.7
$ws.Range("E15").Value = 747.35
$ws.Range("F15").Value = 52
$ws.Range("G15").Value = 737.6

$ws.Range("B16").Value = 1237.55
$ws.Range("C16").Value = 1222
$ws.Range("D16").Value = 1223.3
$ws.Range("E16").Value = 1223
$ws.Range("F16").Value = 84
$ws.Range("G16").Value = 1234.7

$ws.Range("B17").Value = 1371.3
$ws.Range("C17").Value = 1339
$ws.Range("D17").Value = 1347.85
$ws.Range("E17").Value = 1347.6
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 1356.15

$ws.Range("B18").Value = 1947
$ws.Range("C18").Value = 1915.5
$ws.Range("D18").Value = 1935
$ws.Range("E18").Value = 1935.1
$ws.Range("F18").Value = 47
$ws.Range("G18").Value = 1920.65

$ws.Range("B19").Value = 1014
$ws.Range("C19").Value = 996.65
$ws.Range("D19").Value = 1003.9
$ws.Range("E19").Value = 1005.85
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 997.45

$ws.Range("B20").Value = 621.65
$ws.Range("C20").Value = 614.15
$ws.Range("D20").Value = 620.05
$ws.Range("E20").Value = 619.5
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 615.55

$ws.Range("B21").Value = 3188
$ws.Range("C21").Value = 3121
$ws.Range("D21").Value = 3136.8
$ws.Range("E21").Value = 3134.35
$ws.Range("F21").Value = 16
$ws.Range("G21").Value = 3180.75

$ws.Range("B22").Value = 287.05
$ws.Range("C22").Value = 282.6
$ws.Range("D22").Value = 284.4
$ws.Range("E22").Value = 284.05
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 284.45

$ws.Range("B23").Value = 426.65
$ws.Range("C23").Value = 420.3
$ws.Range("D23").Value = 422.5
$ws.Range("E23").Value = 422.5
$ws.Range("F23").Value = 86
$ws.Range("G23").Value = 420.6

$ws.Range("B24").Value = 2766.55
$ws.Range("C24").Value = 2740
$ws.Range("D24").Value = 2749
$ws.Range("E24").Value = 2744.2
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = 2751.95

$ws.Range("B25").Value = 802.7
$ws.Range("C25").Value = 793
$ws.Range("D25").Value = 799
$ws.Range("E25").Value = 799.75
$ws.Range("F25").Value = 143
$ws.Range("G25").Value = 797

$ws.Range("B26").Value = 795.65
$ws.Range("C26").Value = 782
$ws.Range("D26").Value = 782.45
$ws.Range("E26").Value = 783.85
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 787.05

$ws.Range("B27").Value = 1188.7
$ws.Range("C27").Value = 1152.9
$ws.Range("D27").Value = 1183.5
$ws.Range("E27").Value = 1183.1
$ws.Range("F27").Value = 43
$ws.Range("G27").Value = 1161.75

$ws.Range("B28").Value = 937
$ws.Range("C28").Value = 928.5
$ws.Range("D28").Value = 931.4
$ws.Range("E28").Value = 930.7
$ws.Range("F28").Value = 44
$ws.Range("G28").Value = 930.15

$ws.Range("B29").Value = 466.3
$ws.Range("C29").Value = 459.55
$ws.Range("D29").Value = 461.5
$ws.Range("E29").Value = 461.15
$ws.Range("F29").Value = 80
$ws.Range("G29").Value = 462.7

$ws.Range("B30").Value = 163.78
$ws.Range("C30").Value = 160.39
$ws.Range("D30").Value = 160.85
$ws.Range("E30").Value = 160.66
$ws.Range("F30").Value = 532
$ws.Range("G30").Value = 160.7

$ws.Range("B31").Value = 11468.3
$ws.Range("C31").Value = 11377.05
$ws.Range("D31").Value = 11423.2
$ws.Range("E31").Value = 11421.3
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 11391.5
